$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Append two new BOM rows: a fuse holder and the fuse it takes.
$ws.Range("A21").Value = "Fuse Holder"
$ws.Range("B21").Value = 1
$ws.Range("C21").Value = "https://www.digikey.com/en/products/detail/keystone-electronics/3568/2137306"

$ws.Range("A22").Value = "Fuse"
$ws.Range("B22").Value = 1
$ws.Range("C22").Value = "https://www.digikey.com/en/products/detail/eaton-bussmann-electrical-division/BK-ATM-3/264847"

# Move the selection, matching the saved view state in the commit.
$ws.Range("R17").Select()
